$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 updates
$ws.Range("J6").Value = 2.77
$ws.Range("L6").Value = 3.55
$ws.Range("Q6").Value = 1.98
$ws.Range("U6").Value = 1.72
$ws.Range("V6").Value = 1.9
$ws.Range("W6").Value = 7.4
$ws.Range("X6").Value = 10.5
$ws.Range("Y6").Value = 9
$ws.Range("AA6").Value = 18.5
$ws.Range("AB6").Value = 29
$ws.Range("AC6").Value = 9
$ws.Range("AD6").Value = 6.1
$ws.Range("AG6").Value = 9.5
$ws.Range("AN6").Value = 4.1
$ws.Range("AR6").Value = 75
$ws.Range("AZ6").Value = 75
$ws.Range("BA6").Value = 100

# Row 18 updates
$ws.Range("G18").Value = 2.1
$ws.Range("I18").Value = 3.9
$ws.Range("L18").Value = 4.5
$ws.Range("W18").Value = 6
$ws.Range("Y18").Value = 10
$ws.Range("AF18").Value = 67
$ws.Range("AH18").Value = 17
$ws.Range("AY18").Value = 34

# Row 23 updates
$ws.Range("Q23").Value = 1.57
$ws.Range("R23").Value = 2.35
